# IObject.xlsx / "DataNode" sheet: fix duplicate Id column header.
#
# The header row (row 1) has a stray "ID" label in column B that should
# read "Id" (same text/format as column A's "Id" header), matching the
# commit's "fix config's bug". Copy A1's formatting onto B1 so its style
# (border etc.) matches A1 exactly, then correct the label text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B1").Value = "Id"
